$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price text that looks numeric (e.g. "261.14").
# Force text formatting first so Excel does not silently convert these
# assignments into numeric values, then restore the default "Normal" style
# so the cell keeps looking like the rest of the sheet (no explicit style).
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "43.820.45"
$ws.Range("E2").Value = "  +2.20%  "
Set-TextValue $ws.Range("D3") "2.214.86"
$ws.Range("E3").Value = "  +0.34%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "261.14"
$ws.Range("E5").Value = "  +2.29%  "
Set-TextValue $ws.Range("D6") "86.48"
$ws.Range("E6").Value = "  +14.10%  "
Set-TextValue $ws.Range("D7") "0.616"
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  +0.00%  "
Set-TextValue $ws.Range("D9") "0.604"
$ws.Range("E9").Value = "  +2.04%  "
Set-TextValue $ws.Range("D10") "45.16"
$ws.Range("E10").Value = "  +8.67%  "
$ws.Range("E11").Value = "  +1.51%  "
Set-TextValue $ws.Range("D12") "7.44"
$ws.Range("E12").Value = "  +8.35%  "
$ws.Range("E13").Value = "  +1.80%  "
Set-TextValue $ws.Range("D14") "2.548.97"
$ws.Range("E14").Value = "  +0.47%  "
Set-TextValue $ws.Range("D15") "14.48"
$ws.Range("E15").Value = "  +0.81%  "
Set-TextValue $ws.Range("D16") "2.211.95"
$ws.Range("E16").Value = "  +0.21%  "
Set-TextValue $ws.Range("D17") "0.785"
$ws.Range("E17").Value = "  +0.91%  "
Set-TextValue $ws.Range("D18") "43.769.03"
$ws.Range("E18").Value = "  +2.32%  "
$ws.Range("E19").Value = "  +1.36%  "
Set-TextValue $ws.Range("D20") "5.94"
$ws.Range("E20").Value = "  +0.54%  "
Set-TextValue $ws.Range("D21") "69.83"
$ws.Range("E21").Value = "  -1.91%  "
Set-TextValue $ws.Range("D22") "2.36"
$ws.Range("E22").Value = "  +8.23%  "
Set-TextValue $ws.Range("D23") "231.66"
$ws.Range("E23").Value = "  +0.99%  "
Set-TextValue $ws.Range("D24") "8.98"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  +5.47%  "
Set-TextValue $ws.Range("D27") "10.67"
$ws.Range("E27").Value = "  +0.60%  "
Set-TextValue $ws.Range("D28") "40.23"
$ws.Range("E28").Value = "  +2.44%  "
$ws.Range("E29").Value = "  +3.47%  "
$ws.Range("E30").Value = "  +1.82%  "
Set-TextValue $ws.Range("D31") "174.36"
$ws.Range("E31").Value = "  +0.68%  "
$ws.Range("E32").Value = "  +1.39%  "
Set-TextValue $ws.Range("D33") "0.0870"
$ws.Range("E33").Value = "  +3.18%  "
Set-TextValue $ws.Range("D34") "5.44"
$ws.Range("E34").Value = "  +4.96%  "
$ws.Range("E35").Value = "  +1.99%  "
$ws.Range("E36").Value = "  +5.28%  "
$ws.Range("E37").Value = "  +5.38%  "
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  +2.49%  "
Set-TextValue $ws.Range("D40") "2.93"
$ws.Range("E40").Value = "  +6.41%  "
$ws.Range("E41").Value = "  +0.30%  "
Set-TextValue $ws.Range("D42") "63.22"
$ws.Range("E42").Value = "  +6.26%  "
$ws.Range("E43").Value = "  +5.11%  "
$ws.Range("E44").Value = "  +1.95%  "
Set-TextValue $ws.Range("D45") "100.79"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("E46").Value = "  +1.04%  "
Set-TextValue $ws.Range("D47") "0.0979"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("E48").Value = "  +4.83%  "
$ws.Range("E49").Value = "  +2.15%  "
Set-TextValue $ws.Range("D50") "0.442"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("E51").Value = "  +5.62%  "
